$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 264555
$ws.Range("E2").Value = 13573
$ws.Range("F2").Value = 13573
$ws.Range("G2").Value = 12420
$ws.Range("H2").Value = 9174
$ws.Range("I2").Value = 9043
$ws.Range("J2").Value = 131
$ws.Range("K2").Value = 229670
$ws.Range("L2").Value = 111836
$ws.Range("M2").Value = 117834
$ws.Range("N2").Value = 114314
$ws.Range("O2").Value = 3520
$ws.Range("P2").Value = 17891
$ws.Range("Q2").Value = 28645
$ws.Range("R2").Value = -34513
$ws.Range("S2").Value = 4047
$ws.Range("T2").Value = 29825
$ws.Range("U2").Value = -1180
$ws.Range("V2").Value = 42474
$ws.Range("W2").Value = 5.13
$ws.Range("X2").Value = 3.47
$ws.Range("Y2").Value = 8.21
$ws.Range("Z2").Value = 4.11
$ws.Range("AA2").Value = 94.91
$ws.Range("AB2").Value = 542.52
$ws.Range("AC2").Value = 2527
$ws.Range("AD2").Value = 13.32
$ws.Range("AE2").Value = 31948
$ws.Range("AF2").Value = 1.05
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 1.49
$ws.Range("AI2").Value = 19.78
$ws.Range("AJ2").Value = 357815700

# Row 3
$ws.Range("D3").Value = 283839
$ws.Range("E3").Value = 16256
$ws.Range("F3").Value = 16256
$ws.Range("G3").Value = 14340
$ws.Range("H3").Value = 10235
$ws.Range("I3").Value = 9666
$ws.Range("J3").Value = 569
$ws.Range("K3").Value = 225772
$ws.Range("L3").Value = 98722
$ws.Range("M3").Value = 127050
$ws.Range("N3").Value = 121930
$ws.Range("O3").Value = 5120
$ws.Range("P3").Value = 17891
$ws.Range("Q3").Value = 27266
$ws.Range("R3").Value = -27319
$ws.Range("S3").Value = -1745
$ws.Range("T3").Value = 23650
$ws.Range("U3").Value = 3616
$ws.Range("V3").Value = 42242
$ws.Range("W3").Value = 5.73
$ws.Range("X3").Value = 3.61
$ws.Range("Y3").Value = 8.18
$ws.Range("Z3").Value = 4.49
$ws.Range("AA3").Value = 77.7
$ws.Range("AB3").Value = 581.84
$ws.Range("AC3").Value = 2701
$ws.Range("AD3").Value = 9.09
$ws.Range("AE3").Value = 34076
$ws.Range("AF3").Value = 0.72
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 2.04
$ws.Range("AI3").Value = 18.51
$ws.Range("AJ3").Value = 357815700

# Row 4
$ws.Range("D4").Value = 265041
$ws.Range("E4").Value = 13114
$ws.Range("F4").Value = 13114
$ws.Range("G4").Value = 13162
$ws.Range("H4").Value = 9315
$ws.Range("I4").Value = 9067
$ws.Range("J4").Value = 248
$ws.Range("K4").Value = 248843
$ws.Range("L4").Value = 114219
$ws.Range("M4").Value = 134624
$ws.Range("N4").Value = 129560
$ws.Range("O4").Value = 5064
$ws.Range("P4").Value = 17891
$ws.Range("Q4").Value = 36409
$ws.Range("R4").Value = -31892
$ws.Range("S4").Value = 3079
$ws.Range("T4").Value = 37359
$ws.Range("U4").Value = -950
$ws.Range("V4").Value = 47788
$ws.Range("W4").Value = 4.95
$ws.Range("X4").Value = 3.52
$ws.Range("Y4").Value = 7.21
$ws.Range("Z4").Value = 3.92
$ws.Range("AA4").Value = 84.84
$ws.Range("AB4").Value = 629.12
$ws.Range("AC4").Value = 2534
$ws.Range("AD4").Value = 12.41
$ws.Range("AE4").Value = 36209
$ws.Range("AF4").Value = 0.87
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 1.59
$ws.Range("AI4").Value = 19.73
$ws.Range("AJ4").Value = 357815700

# Row 5
$ws.Range("D5").Value = 277902
$ws.Range("E5").Value = 24616
$ws.Range("F5").Value = 24616
$ws.Range("G5").Value = 23326
$ws.Range("H5").Value = 19371
$ws.Range("I5").Value = 18028
$ws.Range("J5").Value = 1343
$ws.Range("K5").Value = 291597
$ws.Range("L5").Value = 141782
$ws.Range("M5").Value = 149815
$ws.Range("N5").Value = 143735
$ws.Range("O5").Value = 6080
$ws.Range("P5").Value = 17891
$ws.Range("Q5").Value = 67642
$ws.Range("R5").Value = -64811
$ws.Range("S5").Value = 8622
$ws.Range("T5").Value = 65924
$ws.Range("U5").Value = 1718
$ws.Range("V5").Value = 56031
$ws.Range("W5").Value = 8.86
$ws.Range("X5").Value = 6.97
$ws.Range("Y5").Value = 13.19
$ws.Range("Z5").Value = 7.17
$ws.Range("AA5").Value = 94.64
$ws.Range("AB5").Value = 719.51
$ws.Range("AC5").Value = 5038
$ws.Range("AD5").Value = 5.93
$ws.Range("AE5").Value = 40170
$ws.Range("AF5").Value = 0.74
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 1.67
$ws.Range("AI5").Value = 9.92
$ws.Range("AJ5").Value = 357815700

# Row 6
$ws.Range("D6").Value = 243366
$ws.Range("E6").Value = 929
$ws.Range("F6").Value = 929
$ws.Range("G6").Value = -914
$ws.Range("H6").Value = -1794
$ws.Range("I6").Value = -2072
$ws.Range("K6").Value = 331757
$ws.Range("L6").Value = 182895
$ws.Range("M6").Value = 148862
$ws.Range("N6").Value = 139792
$ws.Range("P6").Value = 17891
$ws.Range("Q6").Value = 44841
$ws.Range("R6").Value = -76753
$ws.Range("S6").Value = 29529
$ws.Range("T6").Value = 79422
$ws.Range("U6").Value = -34581
$ws.Range("V6").Value = 85588
$ws.Range("W6").Value = 0.38
$ws.Range("X6").Value = -0.74
$ws.Range("Y6").Value = -1.46
$ws.Range("Z6").Value = -0.58
$ws.Range("AA6").Value = 122.86
$ws.Range("AB6").Value = 698.18
$ws.Range("AC6").Value = -579
$ws.Range("AD6").Value = -31.16
$ws.Range("AE6").Value = 39068
$ws.Range("AF6").Value = 0.46
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 357815700
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 232116
$ws.Range("E7").Value = -15472
$ws.Range("G7").Value = -18318
$ws.Range("H7").Value = -16001
$ws.Range("I7").Value = -16054
$ws.Range("K7").Value = 359208
$ws.Range("L7").Value = 223065
$ws.Range("M7").Value = 136143
$ws.Range("N7").Value = 125810
$ws.Range("P7").Value = 18118
$ws.Range("Q7").Value = 28992
$ws.Range("R7").Value = -72973
$ws.Range("S7").Value = 39790
$ws.Range("T7").Value = 72955
$ws.Range("U7").Value = -47225
$ws.Range("W7").Value = -6.67
$ws.Range("X7").Value = -6.89
$ws.Range("Y7").Value = -12.09
$ws.Range("Z7").Value = -4.63
$ws.Range("AA7").Value = 163.85
$ws.Range("AC7").Value = -4487
$ws.Range("AD7").Value = -3.31
$ws.Range("AE7").Value = 35161
$ws.Range("AF7").Value = 0.42
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 238224
$ws.Range("E8").Value = -1544
$ws.Range("G8").Value = -4050
$ws.Range("H8").Value = -3442
$ws.Range("I8").Value = -3634
$ws.Range("K8").Value = 365271
$ws.Range("L8").Value = 229957
$ws.Range("M8").Value = 135314
$ws.Range("N8").Value = 125995
$ws.Range("P8").Value = 18521
$ws.Range("Q8").Value = 44628
$ws.Range("R8").Value = -41003
$ws.Range("S8").Value = 5028
$ws.Range("T8").Value = 37981
$ws.Range("U8").Value = 7538
$ws.Range("W8").Value = -0.65
$ws.Range("X8").Value = -1.44
$ws.Range("Y8").Value = -2.89
$ws.Range("Z8").Value = -0.95
$ws.Range("AA8").Value = 169.94
$ws.Range("AC8").Value = -1015
$ws.Range("AD8").Value = -15.17
$ws.Range("AE8").Value = 35212
$ws.Range("AF8").Value = 0.44
$ws.Range("AG8").Value = 2
$ws.Range("AH8").Value = 0.02
$ws.Range("AI8").Value = -0.23

# Row 9
$ws.Range("D9").Value = 251684
$ws.Range("E9").Value = 6254
$ws.Range("G9").Value = 3624
$ws.Range("H9").Value = 2857
$ws.Range("I9").Value = 2691
$ws.Range("K9").Value = 374282
$ws.Range("L9").Value = 235435
$ws.Range("M9").Value = 138848
$ws.Range("N9").Value = 129683
$ws.Range("P9").Value = 18521
$ws.Range("Q9").Value = 50583
$ws.Range("R9").Value = -41920
$ws.Range("S9").Value = -627
$ws.Range("T9").Value = 39392
$ws.Range("U9").Value = 10683
$ws.Range("W9").Value = 2.48
$ws.Range("X9").Value = 1.14
$ws.Range("Y9").Value = 2.1
$ws.Range("Z9").Value = 0.77
$ws.Range("AA9").Value = 169.56
$ws.Range("AC9").Value = 752
$ws.Range("AD9").Value = 20.47
$ws.Range("AE9").Value = 36243
$ws.Range("AF9").Value = 0.42
$ws.Range("AG9").Value = 28
$ws.Range("AH9").Value = 0.18
$ws.Range("AI9").Value = 3.74
